{"js": "// Remove the stale \"M2Doc version mismatch\" warning block (and its\n// surrounding padding spaces) that was left in the first paragraph,\n// right before the word \"query\". This fixes the test template so it no\n// longer contains a hard-coded version-mismatch marker.\nconst body = context.document.body;\n\nconst warningText =\n  \"    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    \";\n\nconst results = body.search(warningText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\n// Delete the matched range(s) entirely (padding spaces + marker + message).\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Remove the stale \"M2Doc version mismatch\" warning block (and its\n# surrounding padding spaces) that was left in the first paragraph,\n# right before the word \"query\". This fixes the test template so it no\n# longer contains a hard-coded version-mismatch marker.\n$d = $word.ActiveDocument\n\n$warningText = \"    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    \"\n\n$range = $d.Content\n$found = $range.Find.Execute($warningText)\nif ($found) {\n  $range.Text = \"\"\n}\n"}
